$d = $word.ActiveDocument

# Shorten the delay-code labels (DELxxM/H -> xxM/H), leaving DEL0 untouched.
$d.Content.Find.Execute("DEL30M", $true, $false, $false, $false, $false, $true, 1, $false, "30M", 2)
$d.Content.Find.Execute("DEL45M", $true, $false, $false, $false, $false, $true, 1, $false, "45M", 2)
$d.Content.Find.Execute("DEL1H", $true, $false, $false, $false, $false, $true, 1, $false, "1H", 2)
$d.Content.Find.Execute("DEL2H", $true, $false, $false, $false, $false, $true, 1, $false, "2H", 2)
$d.Content.Find.Execute("DEL4H", $true, $false, $false, $false, $false, $true, 1, $false, "4H", 2)
$d.Content.Find.Execute("DEL8H", $true, $false, $false, $false, $false, $true, 1, $false, "8H", 2)
$d.Content.Find.Execute("DEL12H", $true, $false, $false, $false, $false, $true, 1, $false, "12H", 2)
$d.Content.Find.Execute("DEL24H", $true, $false, $false, $false, $false, $true, 1, $false, "24H", 2)

# Update the shared wording in every "Description" cell (all 8 occurrences).
$d.Content.Find.Execute("qui reçoit la demande mais avec une présence requise", $true, $false, $false, $false, $false, $true, 1, $false, "qui gère la ressource mais avec une présence", 2)
